$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AI (35) width: widen from 11.7109375 to match target 12.7109375 (closest achievable on the character-width grid)
$ws.Columns.Item(35).ColumnWidth = 11.83

# Update randomized connectivity matrix values
$ws.Range("B1").Value = 0.76789154492876133
$ws.Range("AF1").Value = 0.87636020514697099
$ws.Range("AJ1").Value = 0.93517788358714626
$ws.Range("AM2").Value = 0.83546852234485414
$ws.Range("BH3").Value = 0.59913086539862848
$ws.Range("L4").Value = 0.95549693879929154
$ws.Range("U4").Value = 0.78804554205654154
$ws.Range("AR4").Value = 0.99729117392443278
$ws.Range("AV4").Value = 0.78707092661877853
$ws.Range("W5").Value = 0.63146025266558681
$ws.Range("AT5").Value = 0.97885725701127968
$ws.Range("AU5").Value = 0.78626506430720045
$ws.Range("J6").Value = 0.88679407855487025
$ws.Range("BM6").Value = 0.82818802907468081
$ws.Range("BN6").Value = 0.9673976968530853
$ws.Range("E7").Value = 0.93343611233412505
$ws.Range("O7").Value = 0.96149381970311565
$ws.Range("Y7").Value = 0.82794367668515734
$ws.Range("M8").Value = 0.96473384915773652
$ws.Range("AQ8").Value = 0.65638150558812702
$ws.Range("D9").Value = 0.93292794378257693
$ws.Range("BJ9").Value = 0.95109526848808901
$ws.Range("O10").Value = 0.80610251047842107
$ws.Range("BN10").Value = 0.91687981874544566
$ws.Range("J11").Value = 0.71482072692417875
$ws.Range("L11").Value = 0.96053096610330935
$ws.Range("J12").Value = 0.87023846057615439
$ws.Range("Z12").Value = 0.74117555649747857
$ws.Range("AB12").Value = 0.89345127182242801
$ws.Range("C13").Value = 0.9044472883094099
$ws.Range("K13").Value = 0.72524257732816011
$ws.Range("AF13").Value = 0.92168297746944794
$ws.Range("C14").Value = 0.71912081598364352
$ws.Range("E14").Value = 0.72012567644956582
$ws.Range("L14").Value = 0.86522795107182571
$ws.Range("R14").Value = 0.83361049781843644
$ws.Range("BA14").Value = 0.95785092885869549
$ws.Range("P15").Value = 0.99634346893862658
$ws.Range("AQ15").Value = 0.93727909263334164
$ws.Range("Q16").Value = 0.81752643185518747
$ws.Range("AC16").Value = 0.87748311345349905
$ws.Range("A17").Value = 0.55771890514294276
$ws.Range("H17").Value = 0.83993727856430733
$ws.Range("AE17").Value = 0.88995630481392929
$ws.Range("AE18").Value = 0.74498277757034648
$ws.Range("AE19").Value = 0.97384634514054524
$ws.Range("BC19").Value = 0.75945803523182065
$ws.Range("BO19").Value = 0.83452888019425164
$ws.Range("AN20").Value = 0.75687905055307192
$ws.Range("BC20").Value = 0.92454942367547854
$ws.Range("V21").Value = 0.84279243550409166
$ws.Range("AC21").Value = 0.95129809149560951
$ws.Range("AP21").Value = 0.68527592020827799
$ws.Range("AA22").Value = 0.76628111054982073
$ws.Range("BK22").Value = 0.85355462100250667
$ws.Range("B23").Value = 0.61715642441596819
$ws.Range("S23").Value = 0.87581745301124192
$ws.Range("AX23").Value = 0.81586568218446764
$ws.Range("BH23").Value = 0.93650253269390582
$ws.Range("BB24").Value = 0.87667757988387596
$ws.Range("AJ25").Value = 0.70656767502460061
$ws.Range("X26").Value = 0.96936926233543286
$ws.Range("Y26").Value = 0.80273407604869773
$ws.Range("AP26").Value = 0.73192188321083562
$ws.Range("BP26").Value = 0.95020472103504572
$ws.Range("G27").Value = 0.93066965591836137
$ws.Range("AL27").Value = 0.89590182191234602
$ws.Range("BG27").Value = 0.82544312334338832
$ws.Range("AV28").Value = 0.70799727304518845
$ws.Range("BF28").Value = 0.85425764838299978
$ws.Range("BH28").Value = 0.71431681251614543
$ws.Range("H29").Value = 0.84161921340581247
$ws.Range("AA29").Value = 0.58180260883622958
$ws.Range("BE30").Value = 0.58852729707873663
$ws.Range("F31").Value = 0.98667115086410595
$ws.Range("BN32").Value = 0.90951043944732635
$ws.Range("E34").Value = 0.57584467180395116
$ws.Range("AY34").Value = 0.99571828688497199
$ws.Range("BN34").Value = 0.99507864890367381
$ws.Range("AI36").Value = 0.54028928592874337
$ws.Range("AO36").Value = 0.84606429068951572
$ws.Range("BD36").Value = 0.66885172492554634
$ws.Range("L37").Value = 0.91586770296169751
$ws.Range("AI37").Value = 0.96758368466406153
$ws.Range("AU37").Value = 0.98792294010080228
$ws.Range("BH37").Value = 0.73481454446693284
$ws.Range("AW38").Value = 0.8759675863344174
$ws.Range("BH38").Value = 0.97556280834362119
$ws.Range("AK39").Value = 0.95407904951655009
$ws.Range("E40").Value = 0.59967495750345357
$ws.Range("AM40").Value = 0.80077612970429435
$ws.Range("AX40").Value = 0.83864011151660922
$ws.Range("N41").Value = 0.78666161871113993
$ws.Range("AT41").Value = 0.96598299804703625
$ws.Range("AZ41").Value = 0.71166048220630262
$ws.Range("D42").Value = 0.87656876767234015
$ws.Range("AV42").Value = 0.86919578182762014
$ws.Range("AM43").Value = 0.91404315879695286
$ws.Range("U44").Value = 0.57810252003289464
$ws.Range("AO44").Value = 0.90631857798393201
$ws.Range("AQ45").Value = 0.74130606066888605
$ws.Range("BP45").Value = 0.82361526897317239
$ws.Range("AR46").Value = 0.78706275221978061
$ws.Range("BB46").Value = 0.98479384392522118
$ws.Range("L47").Value = 0.90511153477355766
$ws.Range("X47").Value = 0.79793844270993275
$ws.Range("AS47").Value = 0.91802838766826467
$ws.Range("AT47").Value = 0.8090633855756304
$ws.Range("AA49").Value = 0.8383986901816789
$ws.Range("X51").Value = 0.94245865863984113
$ws.Range("AG51").Value = 0.93411734532719859
$ws.Range("BA51").Value = 0.95669139956876315
$ws.Range("BI51").Value = 0.65710830310998225
$ws.Range("AX52").Value = 0.84401428787680011
$ws.Range("AZ53").Value = 0.63803788282026641
$ws.Range("BF53").Value = 0.52501917285092814
$ws.Range("BD55").Value = 0.89890876036416023
$ws.Range("Z56").Value = 0.62200746290950193
$ws.Range("W57").Value = 0.7853656270817102
$ws.Range("AL57").Value = 0.8896348726341704
$ws.Range("AL58").Value = 0.88990395621833551
$ws.Range("AE59").Value = 0.88314695428269019
$ws.Range("Q60").Value = 0.99841657479221246
$ws.Range("A61").Value = 0.91651559334887744
$ws.Range("AD61").Value = 0.86452123208413045
$ws.Range("C62").Value = 0.78047909622243705
$ws.Range("BA62").Value = 0.77045729517916306
$ws.Range("BI63").Value = 0.93731502679473433
$ws.Range("AF64").Value = 0.72786074671506107
$ws.Range("T65").Value = 0.53830471322210283
$ws.Range("AG65").Value = 0.95994060885869237
$ws.Range("BC65").Value = 0.67419358531537243
$ws.Range("BK65").Value = 0.78476437613368355
$ws.Range("BO65").Value = 0.74391949918195432
$ws.Range("BL66").Value = 0.73180246903354695
$ws.Range("BP67").Value = 0.70035084932403235
